# ---------------------------------------------------------------
# Actualizacion desde MV -datos- :
#  - corrige celdas existentes (filas 127, 137, 138)
#  - agrega filas nuevas 147-168 (fechas de Agosto 2021)
# ---------------------------------------------------------------
$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correcciones en filas existentes ---
$ws.Cells.Item(127, 2).Value = 565
$ws.Cells.Item(127, 10).Value = 74
$ws.Cells.Item(137, 2).Value = 370
$ws.Cells.Item(137, 9).Value = 58
$ws.Cells.Item(138, 2).Value = 704
$ws.Cells.Item(138, 9).Value = 64

# --- Filas nuevas 147-168 ---
$cell = $ws.Cells.Item(147, 1)
$cell.Formula = '=TEXT("02-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(147, 2).Value = 68
$ws.Cells.Item(147, 3).Value = 0
$ws.Cells.Item(147, 4).Value = 0
$ws.Cells.Item(147, 5).Value = 0
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 16
$ws.Cells.Item(147, 8).Value = 0
$ws.Cells.Item(147, 9).Value = 16
$ws.Cells.Item(147, 10).Value = 35

$cell = $ws.Cells.Item(148, 1)
$cell.Formula = '=TEXT("03-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(148, 2).Value = 196
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 48
$ws.Cells.Item(148, 5).Value = 0
$ws.Cells.Item(148, 6).Value = 39
$ws.Cells.Item(148, 7).Value = 27
$ws.Cells.Item(148, 8).Value = 45
$ws.Cells.Item(148, 9).Value = 8
$ws.Cells.Item(148, 10).Value = 29

$cell = $ws.Cells.Item(149, 1)
$cell.Formula = '=TEXT("04-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(149, 2).Value = 978
$ws.Cells.Item(149, 3).Value = 0
$ws.Cells.Item(149, 4).Value = 369
$ws.Cells.Item(149, 5).Value = 135
$ws.Cells.Item(149, 6).Value = 139
$ws.Cells.Item(149, 7).Value = 147
$ws.Cells.Item(149, 8).Value = 36
$ws.Cells.Item(149, 9).Value = 115
$ws.Cells.Item(149, 10).Value = 37

$cell = $ws.Cells.Item(150, 1)
$cell.Formula = '=TEXT("05-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(150, 2).Value = 621
$ws.Cells.Item(150, 3).Value = 0
$ws.Cells.Item(150, 4).Value = 286
$ws.Cells.Item(150, 5).Value = 0
$ws.Cells.Item(150, 6).Value = 170
$ws.Cells.Item(150, 7).Value = 53
$ws.Cells.Item(150, 8).Value = 73
$ws.Cells.Item(150, 9).Value = 9
$ws.Cells.Item(150, 10).Value = 30

$cell = $ws.Cells.Item(151, 1)
$cell.Formula = '=TEXT("06-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(151, 2).Value = 740
$ws.Cells.Item(151, 3).Value = 153
$ws.Cells.Item(151, 4).Value = 128
$ws.Cells.Item(151, 5).Value = 32
$ws.Cells.Item(151, 6).Value = 48
$ws.Cells.Item(151, 7).Value = 77
$ws.Cells.Item(151, 8).Value = 174
$ws.Cells.Item(151, 9).Value = 61
$ws.Cells.Item(151, 10).Value = 67

$cell = $ws.Cells.Item(152, 1)
$cell.Formula = '=TEXT("09-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(152, 2).Value = 796
$ws.Cells.Item(152, 3).Value = 96
$ws.Cells.Item(152, 4).Value = 0
$ws.Cells.Item(152, 5).Value = 156
$ws.Cells.Item(152, 6).Value = 180
$ws.Cells.Item(152, 7).Value = 79
$ws.Cells.Item(152, 8).Value = 201
$ws.Cells.Item(152, 9).Value = 74
$ws.Cells.Item(152, 10).Value = 10

$cell = $ws.Cells.Item(153, 1)
$cell.Formula = '=TEXT("10-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(153, 2).Value = 837
$ws.Cells.Item(153, 3).Value = 96
$ws.Cells.Item(153, 4).Value = 48
$ws.Cells.Item(153, 5).Value = 210
$ws.Cells.Item(153, 6).Value = 253
$ws.Cells.Item(153, 7).Value = 26
$ws.Cells.Item(153, 8).Value = 118
$ws.Cells.Item(153, 9).Value = 29
$ws.Cells.Item(153, 10).Value = 58

$cell = $ws.Cells.Item(154, 1)
$cell.Formula = '=TEXT("11-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(154, 2).Value = 1452
$ws.Cells.Item(154, 3).Value = 450
$ws.Cells.Item(154, 4).Value = 246
$ws.Cells.Item(154, 5).Value = 156
$ws.Cells.Item(154, 6).Value = 185
$ws.Cells.Item(154, 7).Value = 57
$ws.Cells.Item(154, 8).Value = 183
$ws.Cells.Item(154, 9).Value = 81
$ws.Cells.Item(154, 10).Value = 94

$cell = $ws.Cells.Item(155, 1)
$cell.Formula = '=TEXT("12-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(155, 2).Value = 776
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 134
$ws.Cells.Item(155, 5).Value = 103
$ws.Cells.Item(155, 6).Value = 192
$ws.Cells.Item(155, 7).Value = 106
$ws.Cells.Item(155, 8).Value = 233
$ws.Cells.Item(155, 9).Value = 5
$ws.Cells.Item(155, 10).Value = 4

$cell = $ws.Cells.Item(156, 1)
$cell.Formula = '=TEXT("13-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(156, 2).Value = 457
$ws.Cells.Item(156, 3).Value = 0
$ws.Cells.Item(156, 4).Value = 48
$ws.Cells.Item(156, 5).Value = 53
$ws.Cells.Item(156, 6).Value = 78
$ws.Cells.Item(156, 7).Value = 53
$ws.Cells.Item(156, 8).Value = 105
$ws.Cells.Item(156, 9).Value = 76
$ws.Cells.Item(156, 10).Value = 45

$cell = $ws.Cells.Item(157, 1)
$cell.Formula = '=TEXT("16-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(157, 2).Value = 1049
$ws.Cells.Item(157, 3).Value = 0
$ws.Cells.Item(157, 4).Value = 92
$ws.Cells.Item(157, 5).Value = 32
$ws.Cells.Item(157, 6).Value = 69
$ws.Cells.Item(157, 7).Value = 314
$ws.Cells.Item(157, 8).Value = 456
$ws.Cells.Item(157, 9).Value = 80
$ws.Cells.Item(157, 10).Value = 5

$cell = $ws.Cells.Item(158, 1)
$cell.Formula = '=TEXT("17-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(158, 2).Value = 301
$ws.Cells.Item(158, 3).Value = 96
$ws.Cells.Item(158, 4).Value = 46
$ws.Cells.Item(158, 5).Value = 0
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 16
$ws.Cells.Item(158, 8).Value = 90
$ws.Cells.Item(158, 9).Value = 32
$ws.Cells.Item(158, 10).Value = 21

$cell = $ws.Cells.Item(159, 1)
$cell.Formula = '=TEXT("18-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(159, 2).Value = 1005
$ws.Cells.Item(159, 3).Value = 191
$ws.Cells.Item(159, 4).Value = 46
$ws.Cells.Item(159, 5).Value = 32
$ws.Cells.Item(159, 6).Value = 158
$ws.Cells.Item(159, 7).Value = 42
$ws.Cells.Item(159, 8).Value = 337
$ws.Cells.Item(159, 9).Value = 69
$ws.Cells.Item(159, 10).Value = 130

$cell = $ws.Cells.Item(160, 1)
$cell.Formula = '=TEXT("19-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(160, 2).Value = 429
$ws.Cells.Item(160, 3).Value = 0
$ws.Cells.Item(160, 4).Value = 48
$ws.Cells.Item(160, 5).Value = 64
$ws.Cells.Item(160, 6).Value = 72
$ws.Cells.Item(160, 7).Value = 0
$ws.Cells.Item(160, 8).Value = 96
$ws.Cells.Item(160, 9).Value = 66
$ws.Cells.Item(160, 10).Value = 84

$cell = $ws.Cells.Item(161, 1)
$cell.Formula = '=TEXT("20-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(161, 2).Value = 459
$ws.Cells.Item(161, 3).Value = 0
$ws.Cells.Item(161, 4).Value = 80
$ws.Cells.Item(161, 5).Value = 84
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Value = 0
$ws.Cells.Item(161, 8).Value = 109
$ws.Cells.Item(161, 9).Value = 156
$ws.Cells.Item(161, 10).Value = 29

$cell = $ws.Cells.Item(162, 1)
$cell.Formula = '=TEXT("23-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(162, 2).Value = 575
$ws.Cells.Item(162, 3).Value = 96
$ws.Cells.Item(162, 4).Value = 0
$ws.Cells.Item(162, 5).Value = 0
$ws.Cells.Item(162, 6).Value = 32
$ws.Cells.Item(162, 7).Value = 32
$ws.Cells.Item(162, 8).Value = 289
$ws.Cells.Item(162, 9).Value = 84
$ws.Cells.Item(162, 10).Value = 42

$cell = $ws.Cells.Item(163, 1)
$cell.Formula = '=TEXT("24-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(163, 2).Value = 862
$ws.Cells.Item(163, 3).Value = 307
$ws.Cells.Item(163, 4).Value = 297
$ws.Cells.Item(163, 5).Value = 63
$ws.Cells.Item(163, 6).Value = 48
$ws.Cells.Item(163, 7).Value = 0
$ws.Cells.Item(163, 8).Value = 96
$ws.Cells.Item(163, 9).Value = 5
$ws.Cells.Item(163, 10).Value = 47

$cell = $ws.Cells.Item(164, 1)
$cell.Formula = '=TEXT("25-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(164, 2).Value = 1052
$ws.Cells.Item(164, 3).Value = 153
$ws.Cells.Item(164, 4).Value = 0
$ws.Cells.Item(164, 5).Value = 240
$ws.Cells.Item(164, 6).Value = 316
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 165
$ws.Cells.Item(164, 9).Value = 150
$ws.Cells.Item(164, 10).Value = 29

$cell = $ws.Cells.Item(165, 1)
$cell.Formula = '=TEXT("26-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(165, 2).Value = 1487
$ws.Cells.Item(165, 3).Value = 31
$ws.Cells.Item(165, 4).Value = 156
$ws.Cells.Item(165, 5).Value = 450
$ws.Cells.Item(165, 6).Value = 524
$ws.Cells.Item(165, 7).Value = 126
$ws.Cells.Item(165, 8).Value = 145
$ws.Cells.Item(165, 9).Value = 48
$ws.Cells.Item(165, 10).Value = 6

$cell = $ws.Cells.Item(166, 1)
$cell.Formula = '=TEXT("27-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(166, 2).Value = 1092
$ws.Cells.Item(166, 3).Value = 158
$ws.Cells.Item(166, 4).Value = 160
$ws.Cells.Item(166, 5).Value = 31
$ws.Cells.Item(166, 6).Value = 439
$ws.Cells.Item(166, 7).Value = 210
$ws.Cells.Item(166, 8).Value = 84
$ws.Cells.Item(166, 9).Value = 10
$ws.Cells.Item(166, 10).Value = 0

$cell = $ws.Cells.Item(167, 1)
$cell.Formula = '=TEXT("30-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(167, 2).Value = 631
$ws.Cells.Item(167, 3).Value = 0
$ws.Cells.Item(167, 4).Value = 71
$ws.Cells.Item(167, 5).Value = 0
$ws.Cells.Item(167, 6).Value = 250
$ws.Cells.Item(167, 7).Value = 108
$ws.Cells.Item(167, 8).Value = 158
$ws.Cells.Item(167, 9).Value = 25
$ws.Cells.Item(167, 10).Value = 19

$cell = $ws.Cells.Item(168, 1)
$cell.Formula = '=TEXT("31-08-2021","@")'
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(168, 2).Value = 674
$ws.Cells.Item(168, 3).Value = 0
$ws.Cells.Item(168, 4).Value = 205
$ws.Cells.Item(168, 5).Value = 63
$ws.Cells.Item(168, 6).Value = 142
$ws.Cells.Item(168, 7).Value = 0
$ws.Cells.Item(168, 8).Value = 144
$ws.Cells.Item(168, 9).Value = 64
$ws.Cells.Item(168, 10).Value = 56

$excel.CutCopyMode = 0
